$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the original "_GoBack" bookmark near the top of the document.
#    (Word renumbers the remaining bookmark ids automatically on save.)
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Split the sentence about "situaciones de mantenimiento." into three runs
#    by inserting " y solución de errores" right before the final period.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("facilitar la solución a situaciones de mantenimiento", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$insertPoint = $anchor.Start
$ins = $d.Range($insertPoint, $insertPoint)
$ins.InsertAfter(" y solución de errores")

# Wrap the freshly inserted text with a temporary bookmark: adding a bookmark
# around a range forces Word to split the surrounding run so the wrapped text
# becomes its own run. Deleting the bookmark afterwards leaves the run split.
$newText = $d.Content
$newText.Find.Execute(" y solución de errores", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_TempSplitMantenimiento", $newText)
$d.Bookmarks.Item("_TempSplitMantenimiento").Delete()

# ---------------------------------------------------------------------------
# 3) Add a new "_GoBack" bookmark right after the sentence that ends with
#    "...funciones del sistema." (end of that paragraph, zero-width bookmark).
#    Bookmarks.Add on a truly zero-length Range is mishandled by this runtime
#    (it relocates to doc start), so we temporarily insert a marker
#    character, bookmark across it, then delete the marker while keeping the
#    now zero-width bookmark in place.
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("realizar sus actividades para que el usuario realice una correcta utilización de las funciones del sistema.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$markerPos = $target.Start
$target.InsertAfter("X")

$markerRange = $d.Range($markerPos, $markerPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($markerPos, $markerPos + 1)
$markerRange2.Delete()
